# Weekly update: insert a new price row for "Feria Lagunitas de Puerto Montt - Plátano"
# This pushes all existing rows 663..703 down to 664..704 and fills the freed
# row 663 with the newest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 663; Excel shifts rows 663:703 down to 664:704,
# carrying their formatting/values with them (matches target dimension A1:T704).
$ws.Rows.Item(663).Insert()

# Populate the newly inserted row 663 with this week's record.
$ws.Cells.Item(663, 1).Value = 4
$ws.Cells.Item(663, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(663, 3).Value = "Los Lagos"
$ws.Cells.Item(663, 4).Value = 44931
$ws.Cells.Item(663, 5).Value = 10
$ws.Cells.Item(663, 6).Value = "Fruta"
$ws.Cells.Item(663, 7).Value = 100108
$ws.Cells.Item(663, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(663, 9).Value = 100108006
$ws.Cells.Item(663, 10).Value = "Plátano"
$ws.Cells.Item(663, 11).Value = "Sin especificar"
$ws.Cells.Item(663, 12).Value = "Primera Pintón"
$ws.Cells.Item(663, 13).Value = 400
$ws.Cells.Item(663, 14).Value = 25000
$ws.Cells.Item(663, 15).Value = 25000
$ws.Cells.Item(663, 16).Value = 25000
$ws.Cells.Item(663, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(663, 18).Value = "Ecuador"
$ws.Cells.Item(663, 19).Value = 1250
$ws.Cells.Item(663, 20).Value = 20
